# Updated reference designators on Asset_Cal_Info tab from *GL001 to *GL486
$wb = $excel.ActiveWorkbook
$wsMoorings = $wb.Worksheets.Item("Moorings")
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

$ws.Range("A2").Value = "GS05MOAS-GL486-01-FLORDM000"
$ws.Range("A3").Value = "GS05MOAS-GL486-01-FLORDM000"
$ws.Range("A4").Value = "GS05MOAS-GL486-01-FLORDM000"
$ws.Range("A5").Value = "GS05MOAS-GL486-01-FLORDM000"
$ws.Range("A7").Value = "GS05MOAS-GL486-02-DOSTAM000"
$ws.Range("A9").Value = "GS05MOAS-GL486-04-CTDGVM000"
$ws.Range("A11").Value = "GS05MOAS-GL486-00-ENG000000"

# Update the selected cell on the Asset_Cal_Info sheet, then restore the
# originally active sheet (Moorings) so the workbook's active tab is unchanged.
$ws.Range("A13").Select()
$wsMoorings.Activate()
